# Invalid Login Test was added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "InvalidUser"
$ws.Range("B3").Value = "InvalidPassword"
$ws.Range("C3").Formula = "=FALSE"
